$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.282.32'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '1.883.12'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2832'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07780'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("E12").Value = '  -3.26%  '
$ws.Range("D13").Value = '1.886.01'
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.089'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6760'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.98%  '
$ws.Range("D17").Value = '30.291.96'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '2.128.02'
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.398'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007290'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.183'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.390'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.994'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.371'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09700'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.394'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.485'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.135'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04668'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7067'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.59%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.716'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01873'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("E39").Value = '  +6.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.524'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.971'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8655'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4190'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '982.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.288'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.241'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("E51").Value = '  -4.61%  '
